$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column: force text storage (avoid Excel auto-numeric-conversion of dotted price strings)
# by temporarily setting NumberFormat to Text, then restoring default style so no stray
# style is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '75.419.28'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +7.78%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.663.05'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +8.88%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '187.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +12.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '587.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.534'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.194'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +12.37%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.662.43'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.96%  '

$ws.Range("E11").Value = '  +1.38%  '

$ws.Range("E12").Value = '  +6.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.72'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '75.249.44'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.73%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.152.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +9.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000188'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.37%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.55'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +10.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.676.13'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +9.26%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.15'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +28.38%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +10.27%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '371.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +14.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.81%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.801.94'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.94%  '

$ws.Range("E30").Value = '  +0.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0944'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '518.92'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +13.77%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.69'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.21%  '

$ws.Range("E35").Value = '  +8.04%  '

$ws.Range("E36").Value = '  -0.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.58%  '

$ws.Range("E38").Value = '  +6.21%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.16'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.35'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.38%  '

$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.98'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +13.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '169.44'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +25.55%  '

$ws.Range("E44").Value = '  +11.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.329'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.92%  '

$ws.Range("E47").Value = '  +11.74%  '

$ws.Range("E48").Value = '  +2.89%  '

$ws.Range("E49").Value = '  +16.37%  '

$ws.Range("E50").Value = '  +7.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.532'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.52%  '
